$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet1 -> Priyam Garg)
$ws.Name = "Priyam Garg"

# Data: header row + 5 data rows (13 columns, A..M)
$data = @(
    @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result"),
    @("44th","Sunrisers Hyderabad","Priyam Garg","c †Dhoni b Bravo","7","10","0","0","70.00","Chennai Super Kings","Sharjah","September 30","Super Kings won by 6 wickets (with 2 balls remaining)"),
    @("49th","Sunrisers Hyderabad","Priyam Garg","c Tripathi b Varun","21","31","0","1","67.74","Kolkata Knight Riders","Dubai (DSC)","October 03","KKR won by 6 wickets (with 2 balls remaining)"),
    @("52nd","Sunrisers Hyderabad","Priyam Garg","c de Villiers b Christian","15","11","0","1","136.36","Royal Challengers Bangalore","Abu Dhabi","October 06","Sunrisers won by 4 runs"),
    @("55th","Sunrisers Hyderabad","Priyam Garg","c HH Pandya b Bumrah","29","21","2","1","138.09","Mumbai Indians","Abu Dhabi","October 08","Mumbai won by 42 runs"),
    @("40th","Sunrisers Hyderabad","Priyam Garg","c & b Mustafizur Rahman","0","1","0","0","0.00","Rajasthan Royals","Dubai (DSC)","September 27","Sunrisers won by 7 wickets (with 9 balls remaining)")
)

# The whole used range (A1:M6) must hold these as TEXT values (mirrors the
# source file's numberStoredAsText cells), so format as text before writing.
$fullRange = $ws.Range("A1:M6")
$fullRange.NumberFormat = "@"

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $cell = $ws.Cells.Item($r + 1, $c + 1)
        $cell.Value = $rowValues[$c]
    }
}
